$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the now-empty inline string cells in row 11 (A11:C11)
$ws.Range("A11").Value = $null
$ws.Range("B11").Value = $null
$ws.Range("C11").Value = $null

# Add new row 12 data
# A12:C12 need to exist as empty text cells (matching the blank
# inlineStr-style cells used elsewhere in this sheet). A plain "" assignment
# removes the cell entirely, so force a text-typed blank via the classic
# leading-apostrophe "quote prefix", then strip the resulting formatting.
$ws.Range("A12").Value = "'"
$ws.Range("A12").ClearFormats()
$ws.Range("B12").Value = "'"
$ws.Range("B12").ClearFormats()
$ws.Range("C12").Value = "'"
$ws.Range("C12").ClearFormats()
$ws.Range("D12").Value = 10
$ws.Range("E12").Value = 1
$ws.Range("F12").Value = 10
$ws.Range("G12").Value = "saw"
$ws.Range("H12").NumberFormat = "@"
$ws.Range("H12").Value = "2025-07-08"
$ws.Range("H12").ClearFormats()
